$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# WeeklyEntryMatrix sheet
#   - header "Time Record Period" duplicated to new col I
#   - project value split into three related cells (A2 / G2 / H2)
#   - filter placeholder text updated, old value moved to new col I
# ---------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("WeeklyEntryMatrix")

# New header cell mirroring column F ("Time Record Period")
$wsWeekly.Range("I1").Value = "Time Record Period"

# Old project name moves to G2, keeping the bold style it had in A2
$wsWeekly.Range("G2").Value = "E - Akin_PIMCO-FAS-102657"
$wsWeekly.Range("G2").Font.Bold = $true

# New full engagement name cell
$wsWeekly.Range("H2").Value = "E - Akin_PIMCO-Akin Gump Strauss Hauer & Feld LLP-FVA-102657"
$wsWeekly.Range("H2").Font.Bold = $true

# A2 now holds the shortened project name, still bold
$wsWeekly.Range("A2").Value = "Akin_PIMCO-Akin Gump Strauss Hauer & Feld LLP-FVA-102657"
$wsWeekly.Range("A2").Font.Bold = $true

# Old F2 placeholder value moves to new column I, keeping its text format
$wsWeekly.Range("I2").Value = "(Select Project)"
$wsWeekly.Range("I2").NumberFormat = "@"

# F2 gets the new filter placeholder text
$wsWeekly.Range("F2").Value = "Type to filter projects..."

$wsWeekly.Columns.Item(1).AutoFit()
$wsWeekly.Columns.Item(7).AutoFit()
$wsWeekly.Columns.Item(8).AutoFit()

$wsWeekly.Range("A6").Select()

# ---------------------------------------------------------------
# SummaryLogs sheet
#   - old engagement value moves to new col G
#   - A2 gets shortened project name, now bold
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("SummaryLogs")

$wsSummary.Range("G2").Value = "Engagement Akin_PIMCO-FAS-102657"

$wsSummary.Range("A2").Value = "Akin_PIMCO-Akin Gump Strauss Hauer & Feld LLP-FVA-102657"
$wsSummary.Range("A2").Font.Bold = $true

$wsSummary.Columns.Item(1).AutoFit()

$wsSummary.Range("A2").Select()

# ---------------------------------------------------------------
# DetailLogs sheet
#   - A2 gets shortened project name, now bold (no new column here)
# ---------------------------------------------------------------
$wsDetail = $wb.Worksheets.Item("DetailLogs")

$wsDetail.Range("A2").Value = "Akin_PIMCO-Akin Gump Strauss Hauer & Feld LLP-FVA-102657"
$wsDetail.Range("A2").Font.Bold = $true

$wsDetail.Columns.Item(1).AutoFit()

$wsDetail.Range("G6").Select()

# ---------------------------------------------------------------
# Users sheet - just move the selection (no data change).
# Select this sheet LAST so it stays the active/visible tab, matching
# the original workbook's tabSelected state.
# ---------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("B3").Select()
